$d = $word.ActiveDocument

# Locate the three list paragraphs under "Caja:" that need to be restyled /
# have their bookmark relocated:
#   25: "En el cierre NO avisar sobre sobrante"
#   26: "Al ocurrir el cierre con faltante mostrar alerta"
#   27: "Opción a reajuste disminuir o aumentar ingresando comentario (sin afectar el cierre realizado)"
$startPara = $null
$endPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "En el c*ierre NO avisar sobre sobrante*") {
        $startPara = $p
    }
    if ($t -like "Opci?n a reajuste disminuir o aumentar ingresando comentario*") {
        $endPara = $p
        break
    }
}

if ($startPara -eq $null -or $endPara -eq $null) {
    throw "Could not locate the target paragraphs to edit"
}

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)

$newXml = '<w:p w:rsidR="00DD27BE" w:rsidRDefault="00DD27BE" w:rsidP="008A2333"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="2F5496" w:themeColor="accent5" w:themeShade="BF"/><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="2F5496" w:themeColor="accent5" w:themeShade="BF"/><w:lang w:val="es-419"/></w:rPr><w:t>En el cierre NO avisar sobre sobrante</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="2F5496" w:themeColor="accent5" w:themeShade="BF"/><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="2F5496" w:themeColor="accent5" w:themeShade="BF"/><w:lang w:val="es-419"/></w:rPr><w:t>Al ocurrir el cierre con faltante mostrar alerta</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t>Opción a r</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t>eajuste disminuir o aumentar ingresando comentario (sin afectar el cierre realizado)</w:t></w:r></w:p>'

$rng.InsertXML($newXml) | Out-Null
